$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# Instructions section
$ws.Range("A4").Value = "Instructions"
$ws.Range("A4").Font.Bold = $true
$ws.Range("A5").Value = "Use this file to specify the structure of your model. Any extra sheets will be ignored, so you can include other information in them"

# Page overview section
$ws.Range("A7").Value = "Page Overview"
$ws.Range("A7").Font.Bold = $true

$ws.Range("A8").Value = "Databook Pages"
$ws.Range("B8").Value = "Specify which worksheets will be present in the databook"

$ws.Range("A9").Value = "Compartments"
$ws.Range("B9").Value = "Specify the states that an individual can be in - an individual is only ever in one compartment at a time"

$ws.Range("A10").Value = "Transitions"
$ws.Range("B10").Value = "Specify which transitions between compartments are possible"

$ws.Range("A11").Value = "Characteristics"
$ws.Range("B11").Value = "Specify groups of people (e.g. groups of compartments) for data entry"

$ws.Range("A12").Value = "Parameters"
$ws.Range("B12").Value = "Define how to compute the flows between compartments"

$ws.Range("A13").Value = "Cascades"
$ws.Range("B13").Value = "Use the Cascades sheet to define the cascade if it is more complex than just all characteristics in sequence"

# Update the summary description for the framework (added last so the new shared
# string lands at the end of the table)
$ws.Range("B2").Value = "Framework for a 4-stage cascade model, without vital dynamics or new cases"

# Resize the columns to comfortably fit the new content
$ws.Columns.Item(1).ColumnWidth = 20.666666666666668
$ws.Columns.Item(2).ColumnWidth = 96.33333333333333

# Make the About sheet the active tab/selection
$ws.Range("B5").Select()
$ws.Activate()
